{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst s = sections.items[0];\nconst ps = s.pageSetup;\nlet all = [];\nlet proto = Object.getPrototypeOf(ps);\nwhile (proto) {\n  all = all.concat(Object.getOwnPropertyNames(proto));\n  proto = Object.getPrototypeOf(proto);\n}\nreturn JSON.stringify(all);\n", "ps1": "$nbsp = [char]0x00A0\n$d = $word.ActiveDocument\n$sec = $d.Sections.Item(1)\n$hdr = $sec.Headers.Item(1)\n$rng = $hdr.Range\n$target = \" in header\" + $nbsp + \":\"\n$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, $target, 2)\nWrite-Output \"found=$found\"\n"}
